$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New formula cells inside the existing "K33 block" (rows 33-40) ---
$ws.Range("C37").Formula = "=2*2"
$ws.Range("C39").Formula = "=2.74^12"

# --- New block of rows 56-62 (mirrors the style/layout of rows 46-53) ---

# Row 56: label cell L56 holding the brand-new shared string "all_2" (no special style),
# plus the first data row (SVR).
$ws.Range("L56").Value = "all_2"

$ws.Range("M56").Value = "SVR"
$ws.Range("N56").Value = 0.1588
$ws.Range("O56").Value = 0.05

$ws.Range("M57").Value = "lightgbm"
$ws.Range("N57").Value = 0.1134
$ws.Range("O57").Value = 0.1

$ws.Range("M58").Value = "xgboost"
$ws.Range("N58").Value = 0.1326
$ws.Range("O58").Value = 0.1

$ws.Range("M59").Value = "ridge"
$ws.Range("N59").Value = 0.1101
$ws.Range("O59").Value = 0.2

$ws.Range("M60").Value = "rf"
$ws.Range("N60").Value = 0.1323
$ws.Range("O60").Value = 0.1

$ws.Range("M61").Value = "gbr"
$ws.Range("N61").Value = 0.1113
$ws.Range("O61").Value = 0.1

$ws.Range("O62").Formula = "=SUM(O56:O61)"

# Copy the formatting (font/alignment/number format) from the equivalent rows
# 47:52 (same visual block one section above) onto the new rows 56:61, so the
# new cells pick up the same cell styles already present in the workbook
# (left/vcenter Consolas for the M column, Courier New for the N column) instead
# of creating brand-new style entries.
$ws.Range("M47:O52").Copy()
$ws.Range("M56:O61").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- View state: scroll to show the new rows and select the last used cell ---
$ws.Range("P62").Select()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
